$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 11:39"

# --- Indonesia (row 26) : data refresh only, no reorder ---
Set-RowData 26 190665 3128 136401 46324 0 108 7940

# --- Polonia / Japon swap (rows 47-48) ---
# Polonia moves up to row 47 with refreshed data; Japon drops to row 48 with its prior data.
$ws.Cells.Item(47, 1).Value = "Polonia"
Set-RowData 47 70387 567 52346 15928 0 13 2113
$ws.Cells.Item(48, 1).Value = "Japon"
Set-RowData 48 70268 0 60417 8521 0 0 1330

# --- Austria (row 71) : data refresh only, no reorder ---
Set-RowData 71 29087 358 24828 3524 0 0 735

# --- Croacia (row 90) : data refresh only, no reorder ---
Set-RowData 90 11739 311 8771 2771 0 2 197

# --- Malasia (row 96) : data refresh only, no reorder ---
Set-RowData 96 9391 6 9113 150 0 0 128

# --- Finlandia (row 103) : data refresh only, no reorder ---
Set-RowData 103 8261 36 7350 575 0 0 336

# --- Eslovaquia / Ruanda swap (rows 118-119) ---
$ws.Cells.Item(118, 1).Value = "Eslovaquia"
Set-RowData 118 4526 226 2797 1692 0 0 37
$ws.Cells.Item(119, 1).Value = "Ruanda"
Set-RowData 119 4304 0 2191 2095 0 0 18

# --- Eslovenia moves above Gambia and Sri Lanka (rows 128-130) ---
$ws.Cells.Item(128, 1).Value = "Eslovenia"
Set-RowData 128 3122 43 2483 504 0 1 135
$ws.Cells.Item(129, 1).Value = "Gambia"
Set-RowData 129 3120 0 1295 1726 0 0 99
$ws.Cells.Item(130, 1).Value = "Sri Lanka"
Set-RowData 130 3115 0 2918 185 0 0 12

# --- Lituania (row 132) : data refresh only, no reorder ---
Set-RowData 132 3040 36 1952 1002 0 0 86
